# Desafio final - adiciona novos lancamentos em "detalhamento" e
# recalcula o cruzamento (val_bruto / status) na aba "relatorio".

$wb = $excel.ActiveWorkbook

$wsDet = $wb.Worksheets.Item("detalhamento")
$wsRel = $wb.Worksheets.Item("relatorio")

# --- 1) Novas linhas na aba "detalhamento" (escritorio_vendas, operadora, val_bruto) ---
$newRows = @(
    @(701,"claro",140),
    @(701,"oi",20),
    @(701,"tim",135),
    @(701,"vivo",304),
    @(610,"claro",260),
    @(610,"oi",135),
    @(610,"tim",540),
    @(610,"vivo",1400),
    @(609,"claro",486),
    @(609,"oi",90),
    @(609,"tim",90),
    @(609,"vivo",60),
    @(201,"claro",1585),
    @(201,"oi",634),
    @(201,"tim",520),
    @(201,"vivo",521),
    @(202,"claro",2399),
    @(202,"oi",602),
    @(202,"tim",310),
    @(202,"vivo",447),
    @(206,"claro",45),
    @(206,"oi",40),
    @(206,"vivo",15),
    @(613,"claro",90),
    @(613,"oi",70),
    @(613,"tim",35),
    @(613,"vivo",35),
    @(616,"claro",411),
    @(616,"oi",95),
    @(616,"tim",50),
    @(616,"vivo",70),
    @(205,"claro",123),
    @(401,"claro",500),
    @(401,"oi",375),
    @(401,"tim",125),
    @(401,"vivo",343)
)

# primeira linha livre depois dos dados existentes (linha 127 -> 128)
$detLastRow = $wsDet.Cells.Item(1,1).End(-4121).Row
$r = $detLastRow + 1
foreach ($row in $newRows) {
    $wsDet.Cells.Item($r,1).Value = $row[0]
    $wsDet.Cells.Item($r,2).Value = $row[1]
    $wsDet.Cells.Item($r,3).Value = $row[2]
    $r = $r + 1
}
$detLastRow = $r - 1

# --- 2) Monta tabela de consulta (escritorio_vendas|operadora -> val_bruto) a partir de "detalhamento" ---
$lookup = @{}
for ($i = 2; $i -le $detLastRow; $i++) {
    $key = [string]$wsDet.Cells.Item($i,1).Value2 + "|" + [string]$wsDet.Cells.Item($i,2).Value2
    $lookup[$key] = $wsDet.Cells.Item($i,3).Value2
}

# --- 3) Recalcula val_bruto (coluna D) e status (coluna E) em "relatorio" ---
$relLastRow = $wsRel.Cells.Item(1,1).End(-4121).Row
for ($i = 2; $i -le $relLastRow; $i++) {
    $key = [string]$wsRel.Cells.Item($i,1).Value2 + "|" + [string]$wsRel.Cells.Item($i,2).Value2
    if ($lookup.ContainsKey($key)) {
        $valBruto = $lookup[$key]
    } else {
        $valBruto = 0
    }
    $wsRel.Cells.Item($i,4).Value = $valBruto

    $valLiq = $wsRel.Cells.Item($i,3).Value2
    if ($valLiq -eq $valBruto) {
        $wsRel.Cells.Item($i,5).Value = "ok"
    } else {
        $wsRel.Cells.Item($i,5).Value = "Alerta"
    }
}
